$d = $word.ActiveDocument

# ===========================================================================
# Edit 1: "... considering the object as the derived class)"
#      -> "... considering the object as the parent class)"
# The run is split into three runs: prefix, "parent", " class)".
# ===========================================================================
$ctx1 = $d.Content
$ctx1.Find.Execute("considering the object as the derived class", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$word1 = $ctx1.Duplicate
$word1.Find.Execute("derived", $true, $true, $false, $false, $false, $true, 0, $false, "", 0)
$word1.Text = "parent"
$word1.Bold = 1
$word1.Bold = 0

# ===========================================================================
# Edit 2: "we can't even do the new functions that we defined for the children class!"
#      -> "we can't even use the new functions that we defined for the children class!"
# The run is split into three runs: "we can't even ", "use", " the new ...".
# ===========================================================================
$ctx2 = $d.Content
$ctx2.Find.Execute("we can't even do the new functions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$word2 = $ctx2.Duplicate
$word2.Find.Execute("do", $true, $true, $false, $false, $false, $true, 0, $false, "", 0)
$word2.Text = "use"
$word2.Bold = 1
$word2.Bold = 0

# ===========================================================================
# Edit 3: "only contains virtual functions that maps to the function "
#      -> "only contains virtual functions, which maps to the function "
# The (highlighted) run is split into five runs, all keeping the yellow
# highlight: "only contains virtual functions", ",", " ", "which",
# " maps to the function ".
# ===========================================================================
$ctx3 = $d.Content
$ctx3.Find.Execute("only contains virtual functions that maps to the function", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$thatRange = $ctx3.Duplicate
$thatRange.Find.Execute(" that ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$thatRange.Text = ", which "

$ctx3b = $d.Content
$ctx3b.Find.Execute("only contains virtual functions, which maps to the function", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$commaRange = $ctx3b.Duplicate
$commaRange.Find.Execute(",", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$commaStart = $commaRange.Start
$commaEnd = $commaRange.End

$whichRange = $ctx3b.Duplicate
$whichRange.Find.Execute("which", $true, $true, $false, $false, $false, $true, 0, $false, "", 0)
$whichStart = $whichRange.Start
$whichEnd = $whichRange.End

$rComma = $d.Range($commaStart, $commaEnd)
$rComma.Bold = 1
$rComma.Bold = 0

$rSpace = $d.Range($commaEnd, $whichStart)
$rSpace.Bold = 1
$rSpace.Bold = 0

$rWhich = $d.Range($whichStart, $whichEnd)
$rWhich.Bold = 1
$rWhich.Bold = 0
